$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$url = "https://www.mouser.com/ProductDetail/E-Switch/JN2UEENAGX?qs=%252BZnE%2FxbLNR81sZ2W8kOt9Q%3D%3D"

# Add the hyperlink first (display text defaults to the address); set the
# visible cell text afterward so it doesn't get used as the hyperlink's
# display text.
$ws.Hyperlinks.Add($ws.Range("I2"), $url, [Type]::Missing, [Type]::Missing, $url) | Out-Null

$ws.Range("I2").Value = "JN2UEENAGX E-Switch | Mouser"
$ws.Range("I2").Style = $ws.Range("H1").Style

$ws.Range("I2").Select() | Out-Null
